$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45159,13,100114002,"Camote","Sin especificar","Primera",610,18000,19000,18500,"`$/caja 18 kilos","Perú",1028,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45159,13,100114002,"Camote","Sin especificar","Primera",520,17000,18000,17500,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44845,13,100114002,"Camote","Sin especificar","Primera",600,18000,18000,18000,"`$/malla 18 kilos","Perú",1000,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44270,13,100114002,"Camote","Sin especificar","Primera",16000,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44435,13,100114002,"Camote","Sin especificar","Primera",1150,12000,13000,12500,"`$/malla 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44592,13,100114002,"Camote","Sin especificar","Primera",1000,11000,12000,11600,"`$/malla 18 kilos","Perú",644,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44326,13,100114002,"Camote","Sin especificar","Primera",1600,10000,10000,10000,"`$/malla 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44578,13,100114002,"Camote","Sin especificar","Primera",1240,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44578,13,100114002,"Camote","Sin especificar","Segunda",610,9000,9000,9000,"`$/malla 18 kilos","Perú",500,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44984,13,100114002,"Camote","Sin especificar","Primera",790,18000,19000,18494,"`$/caja 18 kilos","Perú",1027,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44984,13,100114002,"Camote","Sin especificar","Primera",970,15000,16000,15546,"`$/malla 18 kilos","Perú",864,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44921,13,100114002,"Camote","Sin especificar","Primera",790,17000,18000,17506,"`$/caja 18 kilos","Perú",973,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44921,13,100114002,"Camote","Sin especificar","Primera",520,18000,19000,18500,"`$/malla 18 kilos","Perú",1028,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44585,13,100114002,"Camote","Sin especificar","Primera",1330,11000,12000,11500,"`$/malla 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44585,13,100114002,"Camote","Sin especificar","Segunda",610,10000,10000,10000,"`$/malla 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44424,13,100114002,"Camote","Sin especificar","Primera",790,13000,14000,13506,"`$/malla 18 kilos","Perú",750,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44424,13,100114002,"Camote","Sin especificar","Segunda",520,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44781,13,100114002,"Camote","Sin especificar","Primera",430,12000,13000,12500,"`$/caja 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44781,13,100114002,"Camote","Sin especificar","Primera",700,11000,12000,11500,"`$/malla 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44914,13,100114002,"Camote","Sin especificar","Primera",790,17000,18000,17494,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44914,13,100114002,"Camote","Sin especificar","Primera",520,18000,19000,18500,"`$/malla 18 kilos","Perú",1028,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44263,13,100114002,"Camote","Sin especificar","Primera",1600,9000,9000,9000,"`$/malla 18 kilos","Perú",500,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44998,13,100114002,"Camote","Sin especificar","Primera",790,17000,18000,17494,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44998,13,100114002,"Camote","Sin especificar","Primera",880,14000,15000,14500,"`$/malla 18 kilos","Perú",806,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44214,13,100114002,"Camote","Sin especificar","Primera",1900,12000,13000,12526,"`$/malla 18 kilos","Perú",696,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44935,13,100114002,"Camote","Sin especificar","Primera",610,17000,18000,17500,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44935,13,100114002,"Camote","Sin especificar","Primera",520,16000,17000,16500,"`$/malla 18 kilos","Perú",917,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44802,13,100114002,"Camote","Sin especificar","Primera",520,14000,14000,14000,"`$/caja 18 kilos","Perú",778,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44802,13,100114002,"Camote","Sin especificar","Primera",790,14000,14000,14000,"`$/malla 18 kilos","Perú",778,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44305,13,100114002,"Camote","Sin especificar","Primera",1600,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44868,13,100114002,"Camote","Sin especificar","Primera",400,13000,13000,13000,"`$/caja 18 kilos","Perú",722,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44417,13,100114002,"Camote","Sin especificar","Primera",880,14000,15000,14500,"`$/malla 18 kilos","Perú",806,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44417,13,100114002,"Camote","Sin especificar","Segunda",340,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45138,13,100114002,"Camote","Sin especificar","Primera",520,19000,20000,19500,"`$/caja 18 kilos","Perú",1083,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45138,13,100114002,"Camote","Sin especificar","Primera",700,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44249,13,100114002,"Camote","Sin especificar","Primera",1600,10000,10000,10000,"`$/malla 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44809,13,100114002,"Camote","Sin especificar","Primera",900,14000,15000,14667,"`$/caja 18 kilos","Perú",815,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44809,13,100114002,"Camote","Sin especificar","Primera",1300,12000,13000,12538,"`$/malla 18 kilos","Perú",697,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44760,13,100114002,"Camote","Sin especificar","Primera",430,12000,13000,12500,"`$/caja 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44760,13,100114002,"Camote","Sin especificar","Primera",970,9000,10000,9500,"`$/malla 18 kilos","Perú",528,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45124,13,100114002,"Camote","Sin especificar","Primera",520,19000,20000,19500,"`$/caja 18 kilos","Perú",1083,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45124,13,100114002,"Camote","Sin especificar","Primera",305,16000,16000,16000,"`$/malla 18 kilos","Perú",889,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44725,13,100114002,"Camote","Sin especificar","Primera",520,11000,12000,11500,"`$/caja 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44725,13,100114002,"Camote","Sin especificar","Primera",970,9000,10000,9500,"`$/malla 18 kilos","Perú",528,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44956,13,100114002,"Camote","Sin especificar","Primera",610,17000,18000,17500,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44956,13,100114002,"Camote","Sin especificar","Primera",700,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44550,13,100114002,"Camote","Sin especificar","Primera",1060,11000,12000,11500,"`$/malla 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44550,13,100114002,"Camote","Sin especificar","Segunda",430,10000,10000,10000,"`$/malla 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44963,13,100114002,"Camote","Sin especificar","Primera",610,17000,18000,17500,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44963,13,100114002,"Camote","Sin especificar","Primera",520,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44172,13,100114002,"Camote","Sin especificar","Primera",1600,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45005,13,100114002,"Camote","Sin especificar","Primera",790,17000,18000,17494,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45005,13,100114002,"Camote","Sin especificar","Primera",610,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45152,13,100114002,"Camote","Sin especificar","Primera",610,18000,19000,18500,"`$/caja 18 kilos","Perú",1028,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45152,13,100114002,"Camote","Sin especificar","Primera",790,15000,16000,15494,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44795,13,100114002,"Camote","Sin especificar","Primera",1700,13000,14000,13471,"`$/caja 18 kilos","Perú",748,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44795,13,100114002,"Camote","Sin especificar","Primera",1400,13000,14000,13571,"`$/malla 18 kilos","Perú",754,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44669,13,100114002,"Camote","Sin especificar","Primera",1420,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44403,13,100114002,"Camote","Sin especificar","Primera",1330,11000,12000,11500,"`$/caja 15 kilos granel","Perú",767,15,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45012,13,100114002,"Camote","Sin especificar","Primera",520,17000,18000,17500,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45012,13,100114002,"Camote","Sin especificar","Primera",790,10000,11000,10494,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44382,13,100114002,"Camote","Sin especificar","Primera",1510,8000,9000,8500,"`$/malla 18 kilos","Perú",472,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44179,13,100114002,"Camote","Sin especificar","Primera",1500,10000,11000,10600,"`$/malla 18 kilos","Perú",589,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44298,13,100114002,"Camote","Sin especificar","Primera",1600,14000,14000,14000,"`$/malla 18 kilos","Perú",778,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45145,13,100114002,"Camote","Sin especificar","Primera",790,17000,18000,17494,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45145,13,100114002,"Camote","Sin especificar","Primera",880,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45061,13,100114002,"Camote","Sin especificar","Primera",610,19000,20000,19500,"`$/caja 18 kilos","Perú",1083,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45061,13,100114002,"Camote","Sin especificar","Primera",520,16000,17000,16500,"`$/malla 18 kilos","Perú",917,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44522,13,100114002,"Camote","Sin especificar","Primera",1600,11000,12000,11500,"`$/malla 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44376,13,100114002,"Camote","Sin especificar","Primera",520,9000,10000,9500,"`$/malla 18 kilos","Perú",528,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44991,13,100114002,"Camote","Sin especificar","Primera",790,17000,19000,17987,"`$/caja 18 kilos","Perú",999,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44991,13,100114002,"Camote","Sin especificar","Primera",610,14000,15000,14500,"`$/malla 18 kilos","Perú",806,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44466,13,100114002,"Camote","Sin especificar","Primera",1330,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44466,13,100114002,"Camote","Sin especificar","Segunda",790,9000,9000,9000,"`$/malla 18 kilos","Perú",500,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44410,13,100114002,"Camote","Sin especificar","Primera",970,14000,15000,14505,"`$/malla 18 kilos","Perú",806,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44410,13,100114002,"Camote","Sin especificar","Segunda",340,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44949,13,100114002,"Camote","Sin especificar","Primera",520,18000,19000,18500,"`$/caja 18 kilos","Perú",1028,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44949,13,100114002,"Camote","Sin especificar","Primera",610,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44858,13,100114002,"Camote","Sin especificar","Primera",750,17000,18000,17533,"`$/caja 18 kilos","Perú",974,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44858,13,100114002,"Camote","Sin especificar","Primera",1500,13000,14000,13533,"`$/malla 18 kilos","Perú",752,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45082,13,100114002,"Camote","Sin especificar","Primera",790,19000,20000,19494,"`$/caja 18 kilos","Perú",1083,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45082,13,100114002,"Camote","Sin especificar","Primera",430,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44431,13,100114002,"Camote","Sin especificar","Primera",1150,12000,13000,12500,"`$/malla 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44704,13,100114002,"Camote","Sin especificar","Primera",650,10000,10000,10000,"`$/caja 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44704,13,100114002,"Camote","Sin especificar","Primera",1100,7500,8000,7773,"`$/malla 18 kilos","Perú",432,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44284,13,100114002,"Camote","Sin especificar","Primera",1600,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44753,13,100114002,"Camote","Sin especificar","Primera",520,12000,13000,12500,"`$/caja 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44753,13,100114002,"Camote","Sin especificar","Primera",1060,9000,10000,9500,"`$/malla 18 kilos","Perú",528,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44494,13,100114002,"Camote","Sin especificar","Primera",430,17000,17000,17000,"`$/malla 18 kilos","Perú",944,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44340,13,100114002,"Camote","Sin especificar","Primera",1420,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44340,13,100114002,"Camote","Sin especificar","Segunda",970,8000,8000,8000,"`$/malla 18 kilos","Perú",444,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44746,13,100114002,"Camote","Sin especificar","Primera",610,11000,12000,11500,"`$/caja 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44746,13,100114002,"Camote","Sin especificar","Primera",1060,9000,10000,9500,"`$/malla 18 kilos","Perú",528,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44291,13,100114002,"Camote","Sin especificar","Primera",1600,12000,12000,12000,"`$/malla 18 kilos","Perú",667,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45075,13,100114002,"Camote","Sin especificar","Primera",790,18000,19000,18494,"`$/caja 18 kilos","Perú",1027,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45075,13,100114002,"Camote","Sin especificar","Primera",880,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44970,13,100114002,"Camote","Sin especificar","Primera",790,17000,18000,17494,"`$/caja 18 kilos","Perú",972,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44970,13,100114002,"Camote","Sin especificar","Primera",970,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44655,13,100114002,"Camote","Sin especificar","Primera",1600,8000,9000,8500,"`$/malla 18 kilos","Perú",472,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44319,13,100114002,"Camote","Sin especificar","Primera",1510,10000,11000,10500,"`$/malla 18 kilos","Perú",583,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44767,13,100114002,"Camote","Sin especificar","Primera",430,13000,14000,13500,"`$/caja 18 kilos","Perú",750,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44767,13,100114002,"Camote","Sin especificar","Primera",880,9000,10000,9500,"`$/malla 18 kilos","Perú",528,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44459,13,100114002,"Camote","Sin especificar","Primera",1060,12000,13000,12500,"`$/malla 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45117,13,100114002,"Camote","Sin especificar","Primera",430,19000,20000,19500,"`$/caja 18 kilos","Perú",1083,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",45117,13,100114002,"Camote","Sin especificar","Primera",700,15000,16000,15500,"`$/malla 18 kilos","Perú",861,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44893,13,100114002,"Camote","Sin especificar","Primera",520,20000,20000,20000,"`$/caja 18 kilos","Perú",1111,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44893,13,100114002,"Camote","Sin especificar","Primera",610,15000,15000,15000,"`$/malla 18 kilos","Perú",833,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44557,13,100114002,"Camote","Sin especificar","Primera",1600,11000,12000,11500,"`$/malla 18 kilos","Perú",639,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44557,13,100114002,"Camote","Sin especificar","Segunda",790,9000,9000,9000,"`$/malla 18 kilos","Perú",500,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44242,13,100114002,"Camote","Sin especificar","Primera",1600,10000,10000,10000,"`$/malla 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44571,13,100114002,"Camote","Sin especificar","Primera",1060,12000,13000,12500,"`$/malla 18 kilos","Perú",694,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44571,13,100114002,"Camote","Sin especificar","Segunda",520,10000,10000,10000,"`$/malla 18 kilos","Perú",556,18,"Hortaliza"),
    @(9,"Vega Central Mapocho de Santiago","Metropolitana",44711,13,100114002,"Camote","Sin especificar","Primera",970,14000,15000,14500,"`$/malla 18 kilos","Perú",806,18,"Hortaliza"),
)

$rowCount = $data.Count
$colCount = $data[0].Count
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

# Ensure new rows (181, 182) inherit the date number format used by column D
$dateFormat = $ws.Range("D2").NumberFormat
$ws.Range("D181:D182").NumberFormat = $dateFormat

$ws.Range("A70:R182").Value = $arr

Write-Output "Updated rows 70 to 182 ($rowCount rows x $colCount cols)"
